# Updated cryptos list snapshot values.
# Each Value is set with a leading apostrophe ('value) so Excel stores it as
# literal text (matching the workbook's existing inline-string cells) instead
# of silently reinterpreting numeric-looking strings (e.g. "56.90" -> 56.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''36.402.66'
$ws.Range("E2").Value = '''  +0.30%  '
$ws.Range("D3").Value = '''1.940.28'
$ws.Range("E4").Value = '''  -0.11%  '
$ws.Range("D5").Value = '''242.82'
$ws.Range("E5").Value = '''  -0.91%  '
$ws.Range("D6").Value = '''0.613'
$ws.Range("E6").Value = '''  -1.78%  '
$ws.Range("E7").Value = '''  -0.04%  '
$ws.Range("D8").Value = '''56.90'
$ws.Range("E8").Value = '''  -1.32%  '
$ws.Range("D9").Value = '''0.360'
$ws.Range("E9").Value = '''  -3.57%  '
$ws.Range("D10").Value = '''0.0847'
$ws.Range("E10").Value = '''  -2.97%  '
$ws.Range("D12").Value = '''2.225.80'
$ws.Range("E12").Value = '''  -1.61%  '
$ws.Range("B13").Value = '''Avalanche'
$ws.Range("C13").Value = '''https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '''21.26'
$ws.Range("E13").Value = '''  -2.14%  '
$ws.Range("B14").Value = '''Polygon'
$ws.Range("C14").Value = '''https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.808'
$ws.Range("E14").Value = '''  -4.37%  '
$ws.Range("D15").Value = '''13.44'
$ws.Range("E15").Value = '''  -1.37%  '
$ws.Range("E16").Value = '''  -4.64%  '
$ws.Range("D17").Value = '''1.944.87'
$ws.Range("E17").Value = '''  -1.28%  '
$ws.Range("D18").Value = '''36.346.46'
$ws.Range("E18").Value = '''  +0.42%  '
$ws.Range("D19").Value = '''69.18'
$ws.Range("E19").Value = '''  -2.46%  '
$ws.Range("D20").Value = '''0.0₃0859'
$ws.Range("E20").Value = '''  -4.19%  '
$ws.Range("D21").Value = '''228.26'
$ws.Range("E21").Value = '''  -2.30%  '
$ws.Range("E22").Value = '''  -4.18%  '
$ws.Range("E23").Value = '''  -0.07%  '
$ws.Range("D24").Value = '''2.35'
$ws.Range("E24").Value = '''  -6.51%  '
$ws.Range("D25").Value = '''2.28'
$ws.Range("E25").Value = '''  -0.06%  '
$ws.Range("E26").Value = '''  -4.99%  '
$ws.Range("D27").Value = '''161.76'
$ws.Range("E27").Value = '''  -2.66%  '
$ws.Range("D28").Value = '''0.132'
$ws.Range("E28").Value = '''  +4.14%  '
$ws.Range("D29").Value = '''19.42'
$ws.Range("E29").Value = '''  -4.03%  '
$ws.Range("E30").Value = '''  -1.39%  '
$ws.Range("E31").Value = '''  -6.31%  '
$ws.Range("D32").Value = '''4.56'
$ws.Range("E32").Value = '''  -5.20%  '
$ws.Range("D33").Value = '''0.0615'
$ws.Range("E33").Value = '''  -4.58%  '
$ws.Range("E34").Value = '''  -4.88%  '
$ws.Range("D35").Value = '''6.19'
$ws.Range("E35").Value = '''  +2.74%  '
$ws.Range("E36").Value = '''  -0.14%  '
$ws.Range("E37").Value = '''  -1.17%  '
$ws.Range("D38").Value = '''2.16'
$ws.Range("E38").Value = '''  -1.03%  '
$ws.Range("D39").Value = '''3.16'
$ws.Range("E39").Value = '''  +8.64%  '
$ws.Range("D40").Value = '''0.0984'
$ws.Range("E40").Value = '''  +2.55%  '
$ws.Range("E41").Value = '''  -0.05%  '
$ws.Range("D42").Value = '''0.0209'
$ws.Range("E42").Value = '''  -1.40%  '
$ws.Range("D43").Value = '''1.14'
$ws.Range("E43").Value = '''  -4.68%  '
$ws.Range("D44").Value = '''15.80'
$ws.Range("E44").Value = '''  -0.67%  '
$ws.Range("D45").Value = '''1.339.58'
$ws.Range("E45").Value = '''  -1.39%  '
$ws.Range("E46").Value = '''  -5.16%  '
$ws.Range("D47").Value = '''86.49'
$ws.Range("E47").Value = '''  -3.97%  '
$ws.Range("D48").Value = '''7.14'
$ws.Range("E48").Value = '''  -2.84%  '
$ws.Range("E49").Value = '''  +0.25%  '
$ws.Range("D50").Value = '''2.117.01'
$ws.Range("E50").Value = '''  -1.60%  '
$ws.Range("D51").Value = '''43.63'
$ws.Range("E51").Value = '''  -3.19%  '
